$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 values (after edit)
$ws.Range("A4").Value = 112365929
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("Q4").Value = 578487
$ws.Range("R4").Value = 6711585
$ws.Range("S4").Value = 1
$ws.Range("AW4").Value = "Annelie Hilmerby"
$ws.Range("AX4").Value = "Annelie Hilmerby"

# Row 5 values (after edit)
$ws.Range("A5").Value = 112365939
$ws.Range("I5").Value = "6"
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("Q5").Value = 578480
$ws.Range("R5").Value = 6711580
$ws.Range("S5").Value = 20
$ws.Range("AW5").Value = "FREDRIK Månsson"
$ws.Range("AX5").Value = "FREDRIK Månsson"
